# Updated symbol list with GitHub Actions - refresh coin Price/Volume(1h)
# columns (and re-sync a few Coin/Link rows whose ranking order shifted).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D (Price) and E (Volume(1h)) hold numeric-looking text
# (e.g. "287.35", "1.58%") that must stay TEXT, matching the source feed.
# Setting .Value directly on such strings makes Excel auto-coerce them to
# numbers/percentages, so force the text number format first.

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '287.35'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '1.58%'

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '29.55'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '4.36%'

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.084'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '1.32%'

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.06761'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '3.90%'

# Row 6
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '2.07%'

# Row 7
$ws.Range('B7').Value = 'FTXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.401'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '1.37%'

# Row 8
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.9153'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-0.33%'

# Row 9
$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1601'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '4.22%'

# Row 10
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06931'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '4.90%'

# Row 11
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07665'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '1.51%'

# Row 12
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.02926'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '2.82%'

# Row 13
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.09005'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '0.27%'

# Row 14
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.001577'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.50%'

# Row 15
$ws.Range('B15').Value = 'CoinExToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.04479'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '0.95%'

# Row 16
$ws.Range('B16').Value = 'One'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0006471'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '1.36%'

# Row 17
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.006164'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-0.42%'

# Row 18
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.452'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '0.20%'

# Row 19
$ws.Range('B19').Value = 'GateToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.445'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '2.68%'

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.229'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-0.30%'

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.3198'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '0.52%'

# Row 22
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '1.94%'

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.107'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '3.58%'

# Row 24
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '2.44%'

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.001196'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '1.20%'

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.004133'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-6.78%'

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0001200'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '0.06%'

# Row 28
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '-0.07%'

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04264'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '3.58%'

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006829'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '2.05%'

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1242'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '1.13%'

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002231'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '1.89%'

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.01300'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '7.67%'

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005636'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-0.56%'

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.960'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-0.28%'

# Row 47
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-18.62%'
